$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.377411723136902
$ws.Range("B1").Value = 2.150462865829468
$ws.Range("C1").Value = 4.834510803222656
$ws.Range("D1").Value = 3.529049158096313
$ws.Range("E1").Value = 1.247287631034851
